$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Capture the existing (old layout) data before rearranging -------------
# Old layout (row 1 header): A=Day B=Time C=Module Code D=Module Title
#   E=Hours F=Class Type G=Lecturer H=Room I=Block J=Group K=Level L=Course
$lastRow = 10
$old = @{}
for ($r = 2; $r -le $lastRow; $r++) {
    $rowData = @{}
    foreach ($col in @("A","B","C","D","E","F","G","H","I","J","K","L")) {
        $rowData[$col] = $ws.Range("$col$r").Value2
    }
    $old[$r] = $rowData
}

# --- Row 1: collapse header into a single title cell ------------------------
$ws.Range("B1:L1").ClearContents()
$ws.Range("A1").Value = "Herald College Kathmandu"

# --- Rows 2-10: remap columns to the new layout -----------------------------
# New layout: A=Day B=Time C=Hours D=Module Code E=Module Title F=Class Type
#   G=Lecturer H=Group I=Block J=Room   (old K=Level and L=Course are dropped)
for ($r = 2; $r -le $lastRow; $r++) {
    $rowData = $old[$r]
    $ws.Range("A$r").Value = $rowData["A"]
    $ws.Range("B$r").Value = $rowData["B"]
    $ws.Range("C$r").Value = $rowData["E"]
    $ws.Range("D$r").Value = $rowData["C"]
    $ws.Range("E$r").Value = $rowData["D"]
    $ws.Range("F$r").Value = $rowData["F"]
    $ws.Range("G$r").Value = $rowData["G"]
    $ws.Range("H$r").Value = $rowData["J"]
    $ws.Range("I$r").Value = $rowData["I"]
    $ws.Range("J$r").Value = $rowData["H"]
}

# --- Drop the now-unused trailing columns (K:L) -----------------------------
$ws.Range("K1:L$lastRow").ClearContents()
